$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: close off the current last row-group (row 31) with the
# bottom-border style used at the end of every other group (e.g. row 21). ---
$ws.Range("A21:E21").Copy()
$ws.Range("A31:E31").PasteSpecial(-4122)

# --- Step 2: add a new row-group (rows 32-33) for script SCRIPT/G01P03A/um2403.ssb ---

# Row 32 (English / Translated / Converted + filename + line number)
$ws.Cells.Item(32, 3).Value = ' Go for it, Team [team:]!'
$ws.Cells.Item(33, 3).Value = ' Take the Time Gears back to\n[CS:P]Temporal Tower[CR]!'
$ws.Cells.Item(32, 1).Value = 'SCRIPT/G01P03A/um2403.ssb'
$ws.Cells.Item(32, 5).Value = ' Óïìûëï âðåñæä, Ëïíàîäà\n[team:]!'
$ws.Cells.Item(33, 5).Value = ' Âåñîéóå Šåòóåñîé Âñåíåîé â\n[CS:P]Óåíðïñàìûîôý Áàšîý[CR]!'
$ws.Cells.Item(32, 4).Value = ' Только вперёд, Команда\n[team:]!'
$ws.Cells.Item(33, 4).Value = ' Верните Шестерни Времени в\n[CS:P]Темпоральную Башню[CR]!'

$ws.Cells.Item(32, 2).Value = 211
$ws.Cells.Item(33, 2).Value = 214

# Match the formatting used by equivalent rows elsewhere in the sheet
# (row 22: full 5-column entry with a filename; row 8: a continuation
# row without its own filename).
$ws.Range("A22:E22").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122)
$ws.Range("B8:E8").Copy()
$ws.Range("B33:E33").PasteSpecial(-4122)

$ws.Rows.Item(32).RowHeight = 43.2
$ws.Rows.Item(33).RowHeight = 21.6

# --- Step 3: update the view so the new rows are visible / selected ---
$excel.Goto($ws.Range("D35"))
